# Applies the "Add optional courses feedback" edits to the document.
# Each edit is performed with Find/Replace over the whole document content,
# using enough surrounding context to uniquely identify each insertion point.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false,
                               $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement not found for: $old"
    }
}

# 1) "echipa" -> "echipe" (team size clause)
Replace-Text "(in echipa de 3-4)" "(in echipe de 3-4)"

# 2) Add "ca sa iei 10 la aceasta materie" before the sentence-ending period.
#    (kept within the same paragraph -- the following paragraph starts an
#    italic course title "Introducere ..." and Find does not match across
#    paragraph marks; starts right after "nici" so its own spellStart/
#    spellEnd proofing range is left untouched)
Replace-Text "nu consider ca e nevoie." "nu consider ca e nevoie ca sa iei 10 la aceasta materie."

# 3) Add "la ele" after "trebuie mers" (stop before "intelegi" so the
#    spellStart/spellEnd proofing range around it is left untouched)
Replace-Text "si mi se pare ca trebuie mers ca sa " "si mi se pare ca trebuie mers la ele ca sa "

# 4) Add parenthetical remark after "pentru examen." (stop after
#    "rezolvarile" so its own spellStart/spellEnd range stays untouched)
Replace-Text " pentru examen. " " pentru examen (dar am auzit ca unii s-au descurcat si fara sa mearga deloc). "

# 5) Remove "dupa parerea mea" after "multicele cerinte"
Replace-Text "multicele cerinte dupa parerea mea dar tot" "multicele cerinte dar tot"

# 6) Replace "mi se pare ca e nevoie" with "e util"
Replace-Text "laboratorul si mi se pare ca e nevoie sa participi" "laboratorul si e util sa participi"

# 7) Remove "dupa parerea mea" and "(a cate 2 saptamani fiecare)"
Replace-Text "relativ usor de trecut dupa parerea mea. Materia are 5 module (a cate 2 saptamani fiecare) si pentru fiecare" "relativ usor de trecut. Materia are 5 module si pentru fiecare"

# 8) Remove "in timp real" before "in packet" (stop before "packet" so its
#    own spellStart/spellEnd range stays untouched)
Replace-Text " de profi in timp real in " " de profi in "
